$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M4").Value = 2021
$ws.Range("M4").NumberFormat = "General"
$ws.Range("M4").Font.Bold = $true
$ws.Range("M4").Font.Name = "Times New Roman"
$ws.Range("M4").Font.Size = 9
$ws.Range("M4").HorizontalAlignment = -4152
$ws.Range("M4").VerticalAlignment = -4108
